# Fixing issues found #12 + updating gantt chart
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project")

# Period highlight spinner value (linked to N3 via the "period_selected" named range)
$ws.Range("N3").Value = 13

# Try to keep the Spin form-control's reported value in sync with the linked cell.
try {
    $spinner = $ws.Shapes.Item("Spinner 5")
    $spinner.ControlFormat.Value = 13
} catch {
}

# Row 9 - Embed V8 Version 5.5.5 : Actual duration 2 -> 4
$ws.Range("F9").Value = 4

# Row 10 - Implement CommonJs Modules : Actual duration 2 -> 4
$ws.Range("F10").Value = 4

# Row 11 - OpenGL Bindings : Actual start 0 -> 2, Actual duration 0 -> 3, Percent complete 0.35 -> 1
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 1

# Row 12 - OpenCL Bindings : Actual start 0 -> 4, Actual duration 0 -> 3, Percent complete 0.1 -> 1
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 1

# Row 14 - Relevant Reading : Actual duration 0 -> 6, Percent complete 0.85 -> 1
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 1

# Row 15 - Interim Report : Actual start 0 -> 5, Actual duration 0 -> 9, Percent complete 0.5 -> 1
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 1

# Row 16 - Presentation : Actual duration 0 -> 2, Percent complete 0.1 -> 1
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 1

# Row 17 - Basic Demo : Actual duration 0 -> 2, Percent complete 0.1 -> 1
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1

# Row 18 - Debugging & Testing : Actual duration 0 -> 13, Percent complete 0.4 -> 1
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 1

# Update the active selection to match the saved view state
$ws.Range("BJ20").Select()
